# Minor updates to a few slides
#
# 1. presentation.xml: mark the deck as "do not auto-compress pictures"
#    (File > Options > Advanced > Image Size and Quality). This flag is
#    not exposed as a scriptable property on the Presentation/Application
#    object in the PowerPoint object model, so we still attempt the
#    assignment defensively in case the host maps it through, but it is
#    not expected to have an observable effect via COM automation.
# 2. Slide 16 ("Just-In-Time Compiler"): tweak the wording of the first
#    bullet and merge two runs of another bullet into one, dropping the
#    trailing period in favor of a comma.

$p = $ppt.ActivePresentation

try {
    $p.AutoCompressPictures = $false
} catch {
}

$s = $p.Slides.Item(16)
$shape = $s.Shapes.Item("Rectangle 3")
$tr = $shape.TextFrame.TextRange

# --- Edit 1 -----------------------------------------------------------
# The whole run " is a compiler that converts program code into native
# machine code as the program is running." becomes " is a compiler that
# converts program source code into native machine code as the program
# is running." Replacing the full run span (rather than just the changed
# word) keeps this as a single run instead of splitting it into three.
$runOld1 = " is a compiler that converts program code into native machine code as the program is running."
$runNew1 = " is a compiler that converts program source code into native machine code as the program is running."
$start1 = $tr.Text.IndexOf($runOld1) + 1
$tr.Characters($start1, $runOld1.Length).Text = $runNew1

# --- Edit 2 -----------------------------------------------------------
# Two runs "Execution switches to the compiled version once it " +
# "becomes available." are merged into a single run reading
# "Execution switches to the compiled version once it becomes available,"
$targetOld2 = "Execution switches to the compiled version once it becomes available."
$targetNew2 = "Execution switches to the compiled version once it becomes available,"
$start2 = $tr.Text.IndexOf($targetOld2) + 1
$tr.Characters($start2, $targetOld2.Length).Text = $targetNew2
